# Rename worksheets to new spatial-unit identifiers (rerun LU d2c FeatEng
# for FR cities with new spatial units, and dist models).

$wb = $excel.ActiveWorkbook

$renames = @{
    "summ05630095" = "summ53656317"
    "summ05836349" = "summ53886054"
    "summ06040142" = "summ54150809"
    "summ06262646" = "summ54391042"
    "summ06466730" = "summ54634494"
    "summ06687259" = "summ54876124"
    "summ06911828" = "summ55132992"
    "summ07127574" = "summ55381801"
    "summ07506776" = "summ55620759"
}

foreach ($ws in $wb.Worksheets) {
    $oldName = $ws.Name
    if ($renames.ContainsKey($oldName)) {
        $ws.Name = $renames[$oldName]
    }
}
